# update job id to job name in get info api
# Delete the "beginner / no job" row (row 2) from the "jobs" sheet,
# shifting all subsequent job rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jobs")

$ws.Rows.Item(2).Delete()

# Update the selection on the jobs sheet to match the new active cell.
$ws.Range("C14").Select()
